$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = "4260a89d62d231d5dc0c6ac361634926"
$ws.Range("B89").Value = "e87152e4e1851bf4c49bc294e30a4747"
$ws.Range("B99").Value = "1b844d0ff7647576a38ddef541261364"
$ws.Range("B110").Value = "0e0d82d4e819fe84539fda8790233479"
$ws.Range("B136").Value = "7768967d991c633a866b4d5ceb423eb0"
$ws.Range("B154").Value = "32a91f8ae213efc12bec52b2efd9c838"
$ws.Range("B160").Value = "21addbe63baf9031778664dc1908c872"
$ws.Range("B168").Value = "e1fb82d9f179b83e910d00997b76ac1c"
$ws.Range("B169").Value = "1a863b686730fbe8cae77e719eaac063"
$ws.Range("B175").Value = "84112873dbb2c6d20b187742becbd0fd"
$ws.Range("B180").Value = "0a879f3dfbfda86f4e089c6c39181611"
$ws.Range("B213").Value = "9e4b1289d1e75b1181d053a3cf4151b8"
$ws.Range("B222").Value = "911b3ae705610e19631cd514bb71f936"
$ws.Range("B227").Value = "5075d5c6610e4548acbebf191d915c82"
$ws.Range("B228").Value = "a34b05719b14ad4c4ab9c2edbe1a80d7"
$ws.Range("B229").Value = "5122193fee89c48be5b3ac1a03535d59"
$ws.Range("B232").Value = "36f649a3b7651ac52b1344761ee41956"
$ws.Range("B246").Value = "addcbdf982fe3a6fc52c045b518f9dac"
$ws.Range("B279").Value = "a2d7598e285b549a5c9493403e6765d4"
$ws.Range("B281").Value = "5495cb38037e297d4eb77defa7e3dc73"
$ws.Range("B338").Value = "acedbefc8197ac1ca48678cd1130fff7"
$ws.Range("B342").Value = "b1f4ff2e9fb80bad618a2faf81082758"
$ws.Range("B414").Value = "b138ed0d4ad302a67b875946e1802c72"
$ws.Range("B451").Value = "e22a43fda103bdb5c93c1b421a5283ab"
$ws.Range("B467").Value = "93c108a6e9fbf74b96819831bcc7428a"
$ws.Range("B468").Value = "965a7499da59a0349a680c9489d28f69"
$ws.Range("B486").Value = "8a8ea5484628773f105a6844dcdda88e"
$ws.Range("B488").Value = "f52225e2a95c0f31cea88f128314a6f7"
$ws.Range("B525").Value = "d3dc026784724c69f6fec221b4a52adc"
$ws.Range("B526").Value = "b656f9e97bf53832a71f3412b9d582c4"
$ws.Range("B545").Value = "22b82b61937d6fb3c66f33926f213a4e"
$ws.Range("B559").Value = "63fdeceea56c2a5ca6abe64a0d4f2524"
$ws.Range("B578").Value = "4825a72db27fa8e8eb6e66d5c5ce03a4"
$ws.Range("B584").Value = "5ed984f8eb0b051ee955f2040270127d"
$ws.Range("B596").Value = "d593bce827bad7ebd661a55ad706b170"
$ws.Range("B639").Value = "070a12a9d7474a6008726310ce651c92"
$ws.Range("B712").Value = "3350c87a7461671ebab41ea6b500dcf8"
$ws.Range("B715").Value = "f44e5dac00a9c200a5aff2c554453b86"
$ws.Range("B716").Value = "3307567c2f2ad6de937e82dc7d8f11fe"
$ws.Range("B727").Value = "09fb7f76f86704a00a9203cbee4afc2a"
$ws.Range("B745").Value = "f89dc196e5a5fc0a9ddcc67963634c24"
$ws.Range("B754").Value = "704ce11fa59951b7087f65f2cdfd1331"
$ws.Range("B768").Value = "f2940d7e5d7b469038ac04d4c54fc91e"
$ws.Range("B773").Value = "edb6fc22b3a2bf094e0e4d48edb2efd0"
$ws.Range("B780").Value = "a8b3cd5af70366d721d255138ee5ab7a"
$ws.Range("B798").Value = "f01d01adcf2d2c60b413cf1716da18fb"
$ws.Range("B823").Value = "a9f2c216a8ad1ff0db8d4e682aa596d3"
$ws.Range("B827").Value = "bb6022ce7339569742c7321fd58afbb3"
$ws.Range("B831").Value = "cd8968626ca9b2ba70bbee75c334c5cd"
$ws.Range("B837").Value = "b81d3f9d77e61913e8ced6c8f05faef8"
$ws.Range("B839").Value = "5e8db9485fc2c72556604750b10c731b"
$ws.Range("B842").Value = "d216d9c8bb5d663a29571a2cedab1c35"
$ws.Range("B847").Value = "10dabdb298786292523a3b991f934607"
$ws.Range("B866").Value = "bb7c7d197602886af1ce28d88fc7a77b"
$ws.Range("B867").Value = "dd42d95dca18b99b098c9923f94c9db8"
$ws.Range("B881").Value = "ae930a79b34d4121e5878b9c5625bcba"
$ws.Range("B917").Value = "dacacf29fad5085936c128e9c9853864"
$ws.Range("B941").Value = "aac1950523dfaf19462f8ecc460cd0a8"
